# Sprint1Review.docx edit script
# Applies the changes described in the commit diff using Word COM interop.

$d = $word.ActiveDocument

function Get-ParaByText($text) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $ok = $rng.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Could not find paragraph with text: $text"
    }
    return $rng.Paragraphs(1)
}

function Set-ParaSingleRun($para, $newText) {
    $r = $para.Range
    $r.MoveEnd(1, -1)
    $r.Delete()
    $r.InsertAfter($newText)
}

# ---------------------------------------------------------------------------
# 1) Team names
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("William, Daniel, Dylan", $true, $false, $false, $false, $false, $true, 1, $false, "Riley, Trae, Jonah, Anthony", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) "Put names..." sentence: swap "files/documents" -> "documents/files"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Put names of all team members in all files/documents.", $true, $false, $false, $false, $false, $true, 1, $false, "Put names of all team members in all documents/files.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Burndown Chart bullet: "Nice" -> "Looks good"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Nice", $true, $true, $false, $false, $false, $true, 1, $false, "Looks good", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) Remove the old "Requirements:" heading paragraph (and the blank line
#    right before it) from its old spot - it will be re-added later.
# ---------------------------------------------------------------------------
$reqPara = Get-ParaByText("Requirements:")
$blankBeforeReq = $reqPara.Previous()
$reqPara.Range.Delete()
$blankBeforeReq.Range.Delete()

# ---------------------------------------------------------------------------
# 5) "Good. " (two runs) -> "Surprised how little time was required. Is this accurate?"
# ---------------------------------------------------------------------------
$goodPara = Get-ParaByText("Good. ")
Set-ParaSingleRun $goodPara "Surprised how little time was required. Is this accurate?"

# ---------------------------------------------------------------------------
# 6) "UML Diagram:" -> "Requirements:" (re-added here)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("UML Diagram:", $true, $false, $false, $false, $false, $true, 1, $false, "Requirements:", 2) | Out-Null

# ---------------------------------------------------------------------------
# 7) UML Diagram's "Good" bullet -> "Very good." + " "
# ---------------------------------------------------------------------------
$goodPara2 = Get-ParaByText("Good")
$r = $goodPara2.Range
$r.MoveEnd(1, -1)
$r.Delete()
$r.InsertAfter("Very good.")
$r.Collapse(0)
$r.InsertAfter(" ")

# ---------------------------------------------------------------------------
# 8) "Code:" -> "UML Diagram:" (re-added here)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Code:", $true, $true, $false, $false, $false, $true, 1, $false, "UML Diagram:", 2) | Out-Null

# ---------------------------------------------------------------------------
# 9) "All code runs" -> "What does the dotted line represent? "
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("All code runs", $true, $false, $false, $false, $false, $true, 1, $false, "What does the dotted line represent? ", 2) | Out-Null

# ---------------------------------------------------------------------------
# 10) "Very nice GUI for the neighborhood" -> "Looks good"
#     then insert a blank paragraph + new "Code:" heading right after it.
# ---------------------------------------------------------------------------
$guiPara = Get-ParaByText("Very nice GUI for the neighborhood")
Set-ParaSingleRun $guiPara "Looks good"
$guiPara.Range.InsertParagraphAfter()
$blankPara = $guiPara.Next()
$blankPara.Range.InsertParagraphAfter()
$codePara = $blankPara.Next()
$codePara.Range.Text = "Code:"

# ---------------------------------------------------------------------------
# 11) "Some code is documented, some is not" -> "All code runs"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Some code is documented, some is not", $true, $false, $false, $false, $false, $true, 1, $false, "All code runs", 2) | Out-Null

# ---------------------------------------------------------------------------
# 12) "Code implements all parts of sprint goal" -> "Nice " + "graphics display..."
# ---------------------------------------------------------------------------
$implPara = Get-ParaByText("Code implements all parts of sprint goal")
$r = $implPara.Range
$r.MoveEnd(1, -1)
$r.Delete()
$r.InsertAfter("Nice ")
$r.Collapse(0)
$r.InsertAfter("graphics display with text display backup. Good to show client options.")

# ---------------------------------------------------------------------------
# 13) Remove the blank ListParagraph-styled paragraph and old "Team:" heading
#     that used to follow directly.
# ---------------------------------------------------------------------------
$blankList = $implPara.Next()
$oldTeamHeading = $blankList.Next()
# delete the later paragraph first so the earlier reference stays valid
$oldTeamHeading.Range.Delete()
$blankList.Range.Delete()

# ---------------------------------------------------------------------------
# 14) "Team seems to be functional with no apparent communication problems. " (2 runs)
#     -> "All sprint goals seem to be met." then append the new bullet list and
#     "Team:" section underneath it.
# ---------------------------------------------------------------------------
$teamFuncPara = Get-ParaByText("Team seems to be functional with no apparent communication problems.")
Set-ParaSingleRun $teamFuncPara "All sprint goals seem to be met."

$teamFuncPara.Range.InsertParagraphAfter()
$p1 = $teamFuncPara.Next()
$r = $p1.Range
$r.MoveEnd(1, -1)
$r.InsertAfter("Need to document code better.")
$r.Collapse(0)
$r.InsertBreak(6)

$p1.Range.InsertParagraphAfter()
$p2 = $p1.Next()
$p2.Range.Text = "Team:"

$p2.Range.InsertParagraphAfter()
$p3 = $p2.Next()
$p3.Range.Text = "Team seems to be functional."

$p3.Range.InsertParagraphAfter()
$p4 = $p3.Next()
$p4.Range.Text = "No apparent communication problems."

$p4.Range.InsertParagraphAfter()
$p5 = $p4.Next()
$p5.Range.Text = "Larger team (has 4 members) " + [char]0x2026 + " so a bit more expectation from this team."

# ---------------------------------------------------------------------------
# 15) Final assessment paragraph -> two runs
# ---------------------------------------------------------------------------
$finalPara = Get-ParaByText("Excellent first sprint. Talk to each other about process (stop, start, continue) and see if there can be any improvement.")
$r = $finalPara.Range
$r.MoveEnd(1, -1)
$r.Delete()
$r.InsertAfter("All sprint goals accomplished.")
$r.Collapse(0)
$r.InsertAfter(" I am really looking forward to great work from this team!")
